$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value (serial 45202 = 2023-10-03) for every
# data row (rows 2-218). Bump it by one day to 45203 (2023-10-04) for all rows.
$ws.Range("C2:C218").Value = 45203
